# Add a new column "path_on_model3-a" (column L) to the all_runs sheet,
# mapping several run directories to their path on model3-a.
#
# commit message: "add column with mapped paths for model3-a"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column header (L1) -------------------------------------------------
$ws.Range("L1").Value = "path_on_model3-a"

# --- New column values, written in the same order the author would have
#     typed them (first the "current"/NGFr2 rows, then the older NGF_NP10
#     rows) so that the shared-string table is appended in that order.
#     (Paths use literal double backslashes, matching the source file.) -----
$ws.Range("L11").Value = "A:\\Projects\\2035_TM160_NGFr2_NP04_Path1_02"
$ws.Range("L23").Value = "A:\\Projects\\2035_TM160_NGFr2_NP04_Path4_01"
$ws.Range("L14").Value = "G:\\Projects\\2035_TM160_NGF_r2_NoProject_01"
$ws.Range("L19").Value = "B:\\Projects\\2035_TM160_NGF_r2_NoProject_04"
$ws.Range("L26").Value = "F:\\Projects\\2035_TM160_NGFr2_NP04_Path5_01"
$ws.Range("L29").Value = "H:\\Projects\\2035_TM160_NGFr2_NP04_Path6_01"
$ws.Range("L3").Value  = "X:\\Projects\\2035_TM152_NGF_NP10_Path1a_02"
$ws.Range("L4").Value  = "B:\\Projects\\2035_TM152_NGF_NP10_Path1b_02"
$ws.Range("L5").Value  = "A:\\Projects\\2035_TM152_NGF_NP10_Path1x_01"
$ws.Range("L6").Value  = "G:\\Projects\\2035_TM152_NGF_NP10_Path4_02"

# --- Extend the autofilter over the whole used range, including the new
#     column, and keep the _FilterDatabase defined name in sync. ------------
$ws.Range("A1:L29").AutoFilter()

$fd = $wb.Names.Item(1)
$fd.RefersTo = "=all_runs!`$A`$1:`$L`$29"

# --- Put the selection/active cell roughly where the author left it. ------
$ws.Range("C7").Select()
